$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday")

$tuesdayRows = @(
  @("Jackson", "Sherman", "dnc.jacksonsherman.txt", "Tuesday", "speech", "Cision"),
  @("Mikulski", "Barbara", "dnc.mikulski.txt", "Tuesday", "speech", "Cision"),
  @("Lewis", "John", "dnc.lewis.txt", "Tuesday", "speech", "Cision"),
  @("Amaru", "Na'liah", "dnc.amaru.txt", "Tuesday", "speech", "Cision"),
  @("McAuliffe", "Terry", "dnc.mcauliffe.txt", "Tuesday", "speech", "Cision"),
  @("Gabbard", "Tulsi", "dnc.gabbard.txt", "Tuesday", "speech", "Cision"),
  @("Women", "House", "dnc.housewomen.txt", "Tuesday", "speech", "Cision"),
  @("Carter", "Jason", "dnc.carterjason.txt", "Tuesday", "speech", "Cision"),
  @("Schumer", "Chuck", "dnc.schumer.txt", "Tuesday", "speech", "Cision"),
  @("Banks", "Elizabeth", "dnc.banks.txt", "Tuesday", "speech", "Cision"),
  @("Desmond", "haddeus", "dnc.desmond.txt", "Tuesday", "speech", "Cision"),
  @("Moore", "Anton", "dnc.moore.txt", "Tuesday", "speech", "Cision"),
  @("Haubert", "Dynah", "dnc.haubert.txt", "Tuesday", "speech", "Cision"),
  @("Burdick", "Kate", "dnc.burdick.txt", "Tuesday", "speech", "Cision"),
  @("Parson", "Dustin", "dnc.parson.txt", "Tuesday", "speech", "Cision"),
  @("Mellott", "Danielle", "dnc.mellott.txt", "Tuesday", "speech", "Cision"),
  @("Freeman", "Jelani", "dnc.freeman.txt", "Tuesday", "speech", "Cision"),
  @("Feeney", "Paul", "dnc.feeney.txt", "Tuesday", "speech", "Cision"),
  @("Nelson", "Shyla", "dnc.nelsonshyla.txt", "Tuesday", "speech", "Cision"),
  @("Brazile", "Donna", "dnc.brazile.txt", "Tuesday", "speech", "Cision"),
  @("Banks", "David", "dnc.banksdavid.txt", "Tuesday", "speech", "Cision"),
  @("Holder", "Eric", "dnc.holder.txt", "Tuesday", "speech", "Cision"),
  @("Goldwyn", "Tony", "dnc.goldwyn.txt", "Tuesday", "speech", "Cision"),
  @("Movement", "Mothers", "dnc.movement.txt", "Tuesday", "speech", "Cision"),
  @("Richards", "Cecile", "dnc.richards.txt", "Tuesday", "speech", "Cision"),
  @("Benjamin", "Steve", "dnc.benjamin.txt", "Tuesday", "speech", "Cision"),
  @("Boxer", "Barbara", "dnc.boxer.txt", "Tuesday", "speech", "Cision"),
  @("Messing", "Debra", "dnc.messing.txt", "Tuesday", "speech", "Cision"),
  @("Sweeney", "Joe", "dnc.sweeney.txt", "Tuesday", "speech", "Cision"),
  @("Manning", "Lauren", "dnc.manning.txt", "Tuesday", "speech", "Cision"),
  @("McLay", "Cameron", "dnc.mclay.txt", "Tuesday", "speech", "Cision"),
  @("Crowley", "Joe", "dnc.crowley.txt", "Tuesday", "speech", "Cision"),
  @("Alexander", "Erika", "dnc.alexander.txt", "Tuesday", "speech", "Cision"),
  @("Moore", "Ryan", "dnc.mooreryan.txt", "Tuesday", "speech", "Cision"),
  @("Fererra", "America", "dnc.fererra.txt", "Tuesday", "speech", "Cision"),
  @("Dean", "Howard", "dnc.deanhoward.txt", "Tuesday", "speech", "Cision"),
  @("Klobuchar", "Amy", "dnc.klobuchar.txt", "Tuesday", "speech", "Cision"),
  @("Matul", "Ima", "dnc.matul.txt", "Tuesday", "speech", "Cision"),
  @("Campolo", "Tony", "dnc.campolo.txt", "Tuesday", "benediction", "cision"),
)

# Write the new Tuesday rows starting at row 6 (rows 1-5 already existed).
$startRow = 6
for ($i = 0; $i -lt $tuesdayRows.Length; $i++) {
  $r = $startRow + $i
  $rowData = $tuesdayRows[$i]
  for ($j = 0; $j -lt $rowData.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
  }
}

# Row 4's Source value moves from column F to column G (CNN).
$ws.Cells.Item(4, 6).ClearContents()
$ws.Cells.Item(4, 7).Value = "CNN"

# Column C needs to widen to fit the longer names now present.
$ws.Columns.Item(3).ColumnWidth = 21.67

# Tuesday becomes the active sheet/tab.
$ws.Activate()

# Final selection lands on F45, just past the last data row.
$ws.Range("F45").Select()

Write-Output "Tuesday sheet updated with Cision speeches"
